$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A543").Value = 44007
$ws.Range("A543").NumberFormat = $ws.Range("A542").NumberFormat
$ws.Range("B543").Value = 4038.32
$ws.Range("B543").NumberFormat = $ws.Range("B542").NumberFormat
